# Insert a new data row at row 64 (pushing existing rows 64..162 down to 65..163)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 64; this shifts rows 64-162 down to 65-163
# and updates the sheet dimension automatically.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record.
$ws.Range("A64").Value = 5
$ws.Range("B64").Value = "Macroferia Regional de Talca"
$ws.Range("C64").Value = "Maule"
$ws.Range("D64").Value = 44495
$ws.Range("E64").Value = 7
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100108
$ws.Range("H64").Value = "Tropicales y subtropicales"
$ws.Range("I64").Value = 100108005
$ws.Range("J64").Value = "Piña"
$ws.Range("K64").Value = "Caramelo"
$ws.Range("L64").Value = "Tercera"
$ws.Range("M64").Value = 80
$ws.Range("N64").Value = 20000
$ws.Range("O64").Value = 20000
$ws.Range("P64").Value = 20000
$ws.Range("Q64").Value = "`$/caja 16 unidades"
$ws.Range("R64").Value = "Ecuador"
$ws.Range("S64").Value = 1250
$ws.Range("T64").Value = 16
